$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B/C/D columns match the pattern of rows 4-12 (date, scenario csv, building json)
$ws.Cells.Item(13, 2).Value = 44511
$ws.Cells.Item(13, 2).NumberFormat = "d-mmm"
$ws.Cells.Item(13, 3).Value = "Calls_d.csv"
$ws.Cells.Item(13, 4).Value = "B5.json"

$ws.Cells.Item(14, 2).Value = 44511
$ws.Cells.Item(14, 2).NumberFormat = "d-mmm"
$ws.Cells.Item(14, 3).Value = "Calls_d.csv"
$ws.Cells.Item(14, 4).Value = "B5.json"

$ws.Cells.Item(15, 2).Value = 44511
$ws.Cells.Item(15, 2).NumberFormat = "d-mmm"
$ws.Cells.Item(15, 3).Value = "Calls_d.csv"
$ws.Cells.Item(15, 4).Value = "B5.json"

$ws.Cells.Item(16, 2).Value = 44511
$ws.Cells.Item(16, 2).NumberFormat = "d-mmm"
$ws.Cells.Item(16, 3).Value = "Calls_a.csv"
$ws.Cells.Item(16, 4).Value = "B5.json"

# Log text (column A) and algo name (column E), entered in the order the new
# shared-string table was built: A13, A14, E14, E13, A15, E15, A16, E16
$ws.Cells.Item(13, 1).Value = "Total waiting time: 49011.0,  average waiting time per call: 49.011,  unCompleted calls,0,  certificate, -120713707"
$ws.Cells.Item(14, 1).Value = "Total waiting time: 49931.874732,  average waiting time per call: 49.931874732,  unCompleted calls,2,  certificate, -120713699"
$ws.Cells.Item(14, 5).Value = "FlexGreedyAlgoSortWithFixer"
$ws.Cells.Item(13, 5).Value = "FlexGreedyAlgoNoSortWithFixer"
$ws.Cells.Item(15, 1).Value = "Total waiting time: 48702.0,  average waiting time per call: 48.702,  unCompleted calls,0,  certificate, -80037722"
$ws.Cells.Item(15, 5).Value = "FlexGreedyAlgoReserveSortWithFixer"
$ws.Cells.Item(16, 1).Value = "Total waiting time: 1319.0,  average waiting time per call: 13.19,  unCompleted calls,0,  certificate, -260638911"
$ws.Cells.Item(16, 5).Value = "FlexGreedyAlgoReserveSortWithFixer_a"

# Widen columns A and E to fit the new, longer strings
# (Excel rounds ColumnWidth to whole pixels at 7px/char + 5px padding for this
# font, so these inputs are chosen to land exactly on stored widths 95 / 33.)
$ws.Columns.Item(1).ColumnWidth = 94.285714285714285
$ws.Columns.Item(5).ColumnWidth = 32.285714285714285

# Move selection like Excel would after entering the last row of data
$ws.Range("A20").Select()
